$d = $word.ActiveDocument

# Insert the new run's text at the very beginning of the document.
$r = $d.Range(0, 0)
$r.InsertBefore("Ghbdtn /")

# The inserted text now occupies characters [0, 8) ("Ghbdtn /" is 8 chars).
$new = $d.Range(0, 8)
$new.Font.Name = $new.Font.Name  # no-op, keep font as-is
$new.LanguageID = 1033  # en-US
